# Update Betfair Back/Lay odds data for rows 2-12, columns F:AO (36 cols x 11 rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 11,36
$arr[0,0] = 1.57
$arr[0,1] = 1.6
$arr[0,2] = 8.800000000000001
$arr[0,3] = 10
$arr[0,4] = 3.85
$arr[0,5] = 3.95
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 3.75
$arr[0,9] = 1.34
$arr[0,10] = 1.58
$arr[0,11] = 2.68
$arr[0,12] = 1.15
$arr[0,13] = 7
$arr[0,14] = 1.78
$arr[0,15] = 2.06
$arr[0,16] = 1.1
$arr[0,17] = 2.76
$arr[0,18] = 1000
$arr[0,19] = 1000
$arr[0,20] = 1000
$arr[0,21] = 1000
$arr[0,22] = 4.1
$arr[0,23] = 4.9
$arr[0,24] = 16.5
$arr[0,25] = 85
$arr[0,26] = 7
$arr[0,27] = 9.199999999999999
$arr[0,28] = 29
$arr[0,29] = 150
$arr[0,30] = 32
$arr[0,31] = 42
$arr[0,32] = 120
$arr[0,33] = 530
$arr[0,34] = 80
$arr[0,35] = 320
$arr[1,0] = 1.62
$arr[1,1] = 1.66
$arr[1,2] = 6
$arr[1,3] = 6.8
$arr[1,4] = 4.2
$arr[1,5] = 4.5
$arr[1,6] = 0
$arr[1,7] = 0
$arr[1,8] = 8.4
$arr[1,9] = 1.12
$arr[1,10] = 2.76
$arr[1,11] = 1.53
$arr[1,12] = 1.57
$arr[1,13] = 2.64
$arr[1,14] = 1.38
$arr[1,15] = 3.15
$arr[1,16] = 1.17
$arr[1,17] = 2.52
$arr[1,18] = 1000
$arr[1,19] = 1000
$arr[1,20] = 1000
$arr[1,21] = 1000
$arr[1,22] = 9.6
$arr[1,23] = 8.4
$arr[1,24] = 15
$arr[1,25] = 50
$arr[1,26] = 10.5
$arr[1,27] = 8
$arr[1,28] = 13.5
$arr[1,29] = 44
$arr[1,30] = 22
$arr[1,31] = 18.5
$arr[1,32] = 32
$arr[1,33] = 80
$arr[1,34] = 13
$arr[1,35] = 60
$arr[2,0] = 1.45
$arr[2,1] = 1.48
$arr[2,2] = 11
$arr[2,3] = 13
$arr[2,4] = 4
$arr[2,5] = 4.3
$arr[2,6] = 1.48
$arr[2,7] = 1.1
$arr[2,8] = 3
$arr[2,9] = 1.43
$arr[2,10] = 1.67
$arr[2,11] = 2.32
$arr[2,12] = 1.25
$arr[2,13] = 4.6
$arr[2,14] = 2.54
$arr[2,15] = 1.54
$arr[2,16] = 1.08
$arr[2,17] = 3
$arr[2,18] = 11
$arr[2,19] = 26
$arr[2,20] = 140
$arr[2,21] = 1000
$arr[2,22] = 6.2
$arr[2,23] = 10.5
$arr[2,24] = 55
$arr[2,25] = 1000
$arr[2,26] = 7.4
$arr[2,27] = 12
$arr[2,28] = 40
$arr[2,29] = 1000
$arr[2,30] = 12.5
$arr[2,31] = 22
$arr[2,32] = 90
$arr[2,33] = 1000
$arr[2,34] = 12
$arr[2,35] = 1000
$arr[3,0] = 1.9
$arr[3,1] = 1.95
$arr[3,2] = 4.2
$arr[3,3] = 4.5
$arr[3,4] = 4
$arr[3,5] = 4.1
$arr[3,6] = 1.41
$arr[3,7] = 1.06
$arr[3,8] = 4
$arr[3,9] = 1.3
$arr[3,10] = 2.02
$arr[3,11] = 1.92
$arr[3,12] = 1.39
$arr[3,13] = 3.35
$arr[3,14] = 1.8
$arr[3,15] = 1.04
$arr[3,16] = 1.28
$arr[3,17] = 2.04
$arr[3,18] = 16.5
$arr[3,19] = 17
$arr[3,20] = 34
$arr[3,21] = 100
$arr[3,22] = 9.4
$arr[3,23] = 11.5
$arr[3,24] = 17.5
$arr[3,25] = 55
$arr[3,26] = 12
$arr[3,27] = 14
$arr[3,28] = 19
$arr[3,29] = 130
$arr[3,30] = 30
$arr[3,31] = 19.5
$arr[3,32] = 34
$arr[3,33] = 130
$arr[3,34] = 60
$arr[3,35] = 980
$arr[4,0] = 1.54
$arr[4,1] = 1.56
$arr[4,2] = 5.8
$arr[4,3] = 6.4
$arr[4,4] = 4.9
$arr[4,5] = 5.2
$arr[4,6] = 1.27
$arr[4,7] = 1.03
$arr[4,8] = 6.2
$arr[4,9] = 1.16
$arr[4,10] = 2.78
$arr[4,11] = 1.5
$arr[4,12] = 1.75
$arr[4,13] = 2.24
$arr[4,14] = 1.69
$arr[4,15] = 2.28
$arr[4,16] = 1.18
$arr[4,17] = 2.72
$arr[4,18] = 28
$arr[4,19] = 34
$arr[4,20] = 60
$arr[4,21] = 170
$arr[4,22] = 1000
$arr[4,23] = 12
$arr[4,24] = 24
$arr[4,25] = 70
$arr[4,26] = 12
$arr[4,27] = 10
$arr[4,28] = 18.5
$arr[4,29] = 65
$arr[4,30] = 1000
$arr[4,31] = 14.5
$arr[4,32] = 27
$arr[4,33] = 580
$arr[4,34] = 7.6
$arr[4,35] = 55
$arr[5,0] = 1.53
$arr[5,1] = 1.57
$arr[5,2] = 9.4
$arr[5,3] = 11.5
$arr[5,4] = 3.7
$arr[5,5] = 4.2
$arr[5,6] = 1.63
$arr[5,7] = 1.13
$arr[5,8] = 2.48
$arr[5,9] = 1.6
$arr[5,10] = 1.48
$arr[5,11] = 2.84
$arr[5,12] = 1.16
$arr[5,13] = 6
$arr[5,14] = 2.88
$arr[5,15] = 1.43
$arr[5,16] = 1.1
$arr[5,17] = 2.74
$arr[5,18] = 8.6
$arr[5,19] = 21
$arr[5,20] = 95
$arr[5,21] = 1000
$arr[5,22] = 5.1
$arr[5,23] = 11
$arr[5,24] = 50
$arr[5,25] = 400
$arr[5,26] = 6.8
$arr[5,27] = 12.5
$arr[5,28] = 55
$arr[5,29] = 480
$arr[5,30] = 13.5
$arr[5,31] = 26
$arr[5,32] = 100
$arr[5,33] = 1000
$arr[5,34] = 16.5
$arr[5,35] = 970
$arr[6,0] = 1.75
$arr[6,1] = 1.82
$arr[6,2] = 5.1
$arr[6,3] = 6.2
$arr[6,4] = 3.7
$arr[6,5] = 4.1
$arr[6,6] = 1.47
$arr[6,7] = 1.08
$arr[6,8] = 3.4
$arr[6,9] = 1.36
$arr[6,10] = 1.8
$arr[6,11] = 2.06
$arr[6,12] = 1.29
$arr[6,13] = 3.8
$arr[6,14] = 1.94
$arr[6,15] = 1.92
$arr[6,16] = 1.2
$arr[6,17] = 2.18
$arr[6,18] = 24
$arr[6,19] = 18
$arr[6,20] = 1000
$arr[6,21] = 1000
$arr[6,22] = 8.4
$arr[6,23] = 10.5
$arr[6,24] = 1000
$arr[6,25] = 1000
$arr[6,26] = 18
$arr[6,27] = 21
$arr[6,28] = 65
$arr[6,29] = 1000
$arr[6,30] = 900
$arr[6,31] = 1000
$arr[6,32] = 150
$arr[6,33] = 1000
$arr[6,34] = 85
$arr[6,35] = 1000
$arr[7,0] = 1.67
$arr[7,1] = 1.7
$arr[7,2] = 6.6
$arr[7,3] = 7.6
$arr[7,4] = 3.75
$arr[7,5] = 4.1
$arr[7,6] = 1.43
$arr[7,7] = 1.08
$arr[7,8] = 3.35
$arr[7,9] = 1.36
$arr[7,10] = 1.8
$arr[7,11] = 2.08
$arr[7,12] = 1.3
$arr[7,13] = 3.85
$arr[7,14] = 2.1
$arr[7,15] = 1.76
$arr[7,16] = 1.16
$arr[7,17] = 2.42
$arr[7,18] = 12.5
$arr[7,19] = 20
$arr[7,20] = 55
$arr[7,21] = 1000
$arr[7,22] = 7.4
$arr[7,23] = 9
$arr[7,24] = 990
$arr[7,25] = 120
$arr[7,26] = 8.800000000000001
$arr[7,27] = 10
$arr[7,28] = 42
$arr[7,29] = 130
$arr[7,30] = 15.5
$arr[7,31] = 44
$arr[7,32] = 46
$arr[7,33] = 1000
$arr[7,34] = 12.5
$arr[7,35] = 190
$arr[8,0] = 2.18
$arr[8,1] = 2.34
$arr[8,2] = 3.45
$arr[8,3] = 3.9
$arr[8,4] = 3.25
$arr[8,5] = 3.65
$arr[8,6] = 1.42
$arr[8,7] = 1.07
$arr[8,8] = 3.5
$arr[8,9] = 1.33
$arr[8,10] = 1.87
$arr[8,11] = 1.99
$arr[8,12] = 1.33
$arr[8,13] = 3.5
$arr[8,14] = 1.78
$arr[8,15] = 2.02
$arr[8,16] = 1.35
$arr[8,17] = 1.74
$arr[8,18] = 15
$arr[8,19] = 14.5
$arr[8,20] = 38
$arr[8,21] = 440
$arr[8,22] = 9.800000000000001
$arr[8,23] = 15
$arr[8,24] = 16.5
$arr[8,25] = 200
$arr[8,26] = 980
$arr[8,27] = 40
$arr[8,28] = 19
$arr[8,29] = 330
$arr[8,30] = 48
$arr[8,31] = 65
$arr[8,32] = 65
$arr[8,33] = 390
$arr[8,34] = 55
$arr[8,35] = 80
$arr[9,0] = 2.42
$arr[9,1] = 2.48
$arr[9,2] = 3.35
$arr[9,3] = 3.55
$arr[9,4] = 3.3
$arr[9,5] = 3.35
$arr[9,6] = 1.5
$arr[9,7] = 1.09
$arr[9,8] = 3.25
$arr[9,9] = 1.42
$arr[9,10] = 1.78
$arr[9,11] = 2.26
$arr[9,12] = 1.27
$arr[9,13] = 4.4
$arr[9,14] = 1.9
$arr[9,15] = 2
$arr[9,16] = 1.39
$arr[9,17] = 1.68
$arr[9,18] = 11
$arr[9,19] = 11.5
$arr[9,20] = 23
$arr[9,21] = 65
$arr[9,22] = 9
$arr[9,23] = 7.2
$arr[9,24] = 14.5
$arr[9,25] = 46
$arr[9,26] = 14.5
$arr[9,27] = 12
$arr[9,28] = 19
$arr[9,29] = 60
$arr[9,30] = 34
$arr[9,31] = 29
$arr[9,32] = 46
$arr[9,33] = 120
$arr[9,34] = 26
$arr[9,35] = 75
$arr[10,0] = 2.06
$arr[10,1] = 2.1
$arr[10,2] = 4.2
$arr[10,3] = 4.6
$arr[10,4] = 3.45
$arr[10,5] = 3.55
$arr[10,6] = 1.46
$arr[10,7] = 1.09
$arr[10,8] = 3.4
$arr[10,9] = 1.38
$arr[10,10] = 1.79
$arr[10,11] = 2.16
$arr[10,12] = 1.3
$arr[10,13] = 4
$arr[10,14] = 1.89
$arr[10,15] = 2
$arr[10,16] = 1.29
$arr[10,17] = 1.9
$arr[10,18] = 11.5
$arr[10,19] = 14.5
$arr[10,20] = 32
$arr[10,21] = 100
$arr[10,22] = 8.6
$arr[10,23] = 7.8
$arr[10,24] = 17.5
$arr[10,25] = 60
$arr[10,26] = 12
$arr[10,27] = 10.5
$arr[10,28] = 38
$arr[10,29] = 75
$arr[10,30] = 26
$arr[10,31] = 24
$arr[10,32] = 44
$arr[10,33] = 130
$arr[10,34] = 18.5
$arr[10,35] = 80

$ws.Range("F2:AO12").Value = $arr
